$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D7","D9","D11","D12","D13","D17","D18","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D34","D35","D36","D38","D39","D42","D44","D47","D49")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D2').Value = '72.984.34'
$ws.Range('E2').Value = '  +1.47%  '
$ws.Range('D3').Value = '3.991.65'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '615.38'
$ws.Range('E5').Value = '  +15.30%  '
$ws.Range('D6').Value = '165.80'
$ws.Range('E6').Value = '  +10.25%  '
$ws.Range('D7').Value = '0.687'
$ws.Range('E7').Value = '  -1.02%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '0.759'
$ws.Range('E9').Value = '  +0.74%  '
$ws.Range('E10').Value = '  -1.96%  '
$ws.Range('D11').Value = '58.23'
$ws.Range('E11').Value = '  +7.42%  '
$ws.Range('D12').Value = '0.0000317'
$ws.Range('E12').Value = '  -2.34%  '
$ws.Range('D13').Value = '11.22'
$ws.Range('E13').Value = '  +3.93%  '
$ws.Range('D14').Value = '4.623.24'
$ws.Range('E14').Value = '  -0.73%  '
$ws.Range('D15').Value = '4.004.87'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('E16').Value = '  +6.74%  '
$ws.Range('D17').Value = '14.32'
$ws.Range('E17').Value = '  +1.20%  '
$ws.Range('D18').Value = '20.72'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D20').Value = '72.795.16'
$ws.Range('E20').Value = '  +1.19%  '
$ws.Range('D21').Value = '441.93'
$ws.Range('E21').Value = '  +2.13%  '
$ws.Range('D22').Value = '4.92'
$ws.Range('E22').Value = '  +16.72%  '
$ws.Range('D23').Value = '96.63'
$ws.Range('E23').Value = '  -1.75%  '
$ws.Range('E24').Value = '  -5.45%  '
$ws.Range('D25').Value = '14.63'
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('D26').Value = '4.20'
$ws.Range('E26').Value = '  -3.62%  '
$ws.Range('D27').Value = '11.38'
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('D28').Value = '10.57'
$ws.Range('E28').Value = '  -2.17%  '
$ws.Range('D29').Value = '5.96'
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').Value = '36.24'
$ws.Range('E30').Value = '  -1.94%  '
$ws.Range('D31').Value = '7.80'
$ws.Range('E31').Value = '  -5.52%  '
$ws.Range('D32').Value = '13.93'
$ws.Range('E32').Value = '  +2.59%  '
$ws.Range('E33').Value = '  -3.32%  '
$ws.Range('D34').Value = '49.22'
$ws.Range('E34').Value = '  -2.14%  '
$ws.Range('D35').Value = '72.12'
$ws.Range('E35').Value = '  +6.66%  '
$ws.Range('D36').Value = '633.75'
$ws.Range('E36').Value = '  -6.17%  '
$ws.Range('D37').Value = '0.0₃0909'
$ws.Range('E37').Value = '  +9.59%  '
$ws.Range('D38').Value = '0.436'
$ws.Range('E38').Value = '  -3.83%  '
$ws.Range('D39').Value = '3.53'
$ws.Range('E39').Value = '  +4.83%  '
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('D42').Value = '3.35'
$ws.Range('E42').Value = '  -1.84%  '
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('D44').Value = '10.85'
$ws.Range('E44').Value = '  -1.62%  '
$ws.Range('E45').Value = '  -1.04%  '
$ws.Range('E46').Value = '  -0.19%  '
$ws.Range('D47').Value = '2.66'
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('E48').Value = '  +35.32%  '
$ws.Range('D49').Value = '3.40'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').Value = '2.932.55'
$ws.Range('E50').Value = '  +3.00%  '
$ws.Range('E51').Value = '  -0.53%  '
